# F3SPF_LOG SHEET_BBWM_180523.xlsx
# "Processed EEMs and Saved Pictures" / "Processed EEMs for Poster and took
# pictures. 180530"
#
# Append the 24 new "sjer" (San Joaquin Experimental Range) sample rows
# (rows 110-133) to Sheet1, following the same layout used for every other
# sample group already in the log (e.g. "konz", "harv", "ornl"): four
# dilutions (1x/5x/10x/20x) x six blank/correction combinations each.
#
# Cells are written in the same order the strings were originally typed in
# (column A down, then the "...comb" combinations for row block 1, then the
# remaining A-fill, then column B, then column C c1-template-first, then the
# rest of C, finishing with the numeric D column) so the shared-string table
# ends up built in the same sequence as the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Column A (Name of Raw EEM) -------------------------------------------
$ws.Cells.Item(110, 1).Value = "sjer"

# ---- Column B (Blank) for the first block (rows 110-115) ------------------
$ws.Cells.Item(110, 2).Value = "smplhldr"
$ws.Cells.Item(111, 2).Value = "naco3"
$ws.Cells.Item(112, 2).Value = "sjercomb"
$ws.Cells.Item(114, 2).Value = "sjercomb, naco3"
$ws.Cells.Item(113, 2).Value = "sjercomb, smplhldr"
$ws.Cells.Item(115, 2).Value = "sjercomb, naco3, smplhldr"

# ---- Column A for the other three dilution blocks --------------------------
$ws.Cells.Item(116, 1).Value = "sjer5x"
$ws.Cells.Item(122, 1).Value = "sjer10x"
$ws.Cells.Item(128, 1).Value = "sjer20x"

# ---- Fill in the rest of column A (same value repeated down each block) ---
$ws.Cells.Item(111, 1).Value = "sjer"
$ws.Cells.Item(112, 1).Value = "sjer"
$ws.Cells.Item(113, 1).Value = "sjer"
$ws.Cells.Item(114, 1).Value = "sjer"
$ws.Cells.Item(115, 1).Value = "sjer"

$ws.Cells.Item(117, 1).Value = "sjer5x"
$ws.Cells.Item(118, 1).Value = "sjer5x"
$ws.Cells.Item(119, 1).Value = "sjer5x"
$ws.Cells.Item(120, 1).Value = "sjer5x"
$ws.Cells.Item(121, 1).Value = "sjer5x"

$ws.Cells.Item(123, 1).Value = "sjer10x"
$ws.Cells.Item(124, 1).Value = "sjer10x"
$ws.Cells.Item(125, 1).Value = "sjer10x"
$ws.Cells.Item(126, 1).Value = "sjer10x"
$ws.Cells.Item(127, 1).Value = "sjer10x"

$ws.Cells.Item(129, 1).Value = "sjer20x"
$ws.Cells.Item(130, 1).Value = "sjer20x"
$ws.Cells.Item(131, 1).Value = "sjer20x"
$ws.Cells.Item(132, 1).Value = "sjer20x"
$ws.Cells.Item(133, 1).Value = "sjer20x"

# ---- Column B for the remaining three blocks -------------------------------
$ws.Cells.Item(116, 2).Value = "smplhldr"
$ws.Cells.Item(117, 2).Value = "naco3"
$ws.Cells.Item(118, 2).Value = "sjercomb"
$ws.Cells.Item(119, 2).Value = "sjercomb, smplhldr"
$ws.Cells.Item(120, 2).Value = "sjercomb, naco3"
$ws.Cells.Item(121, 2).Value = "sjercomb, naco3, smplhldr"

$ws.Cells.Item(122, 2).Value = "smplhldr"
$ws.Cells.Item(123, 2).Value = "naco3"
$ws.Cells.Item(124, 2).Value = "sjercomb"
$ws.Cells.Item(125, 2).Value = "sjercomb, smplhldr"
$ws.Cells.Item(126, 2).Value = "sjercomb, naco3"
$ws.Cells.Item(127, 2).Value = "sjercomb, naco3, smplhldr"

$ws.Cells.Item(128, 2).Value = "smplhldr"
$ws.Cells.Item(129, 2).Value = "naco3"
$ws.Cells.Item(130, 2).Value = "sjercomb"
$ws.Cells.Item(131, 2).Value = "sjercomb, smplhldr"
$ws.Cells.Item(132, 2).Value = "sjercomb, naco3"
$ws.Cells.Item(133, 2).Value = "sjercomb, naco3, smplhldr"

# ---- Column C (Corrected Name): the "c1" template row of each block first -
$ws.Cells.Item(110, 3).Value = "sjerc1"
$ws.Cells.Item(116, 3).Value = "sjer5xc1"
$ws.Cells.Item(122, 3).Value = "sjer10xc1"
$ws.Cells.Item(128, 3).Value = "sjer20xc1"

# ---- then the rest of the "sjer" (1x) block --------------------------------
$ws.Cells.Item(111, 3).Value = "sjerc2"
$ws.Cells.Item(112, 3).Value = "sjerc3"
$ws.Cells.Item(113, 3).Value = "sjerc4"
$ws.Cells.Item(114, 3).Value = "sjerc5"
$ws.Cells.Item(115, 3).Value = "sjerc6"

# ---- then the rest of the "sjer5x" block -----------------------------------
$ws.Cells.Item(117, 3).Value = "sjer5xc2"
$ws.Cells.Item(118, 3).Value = "sjer5xc3"
$ws.Cells.Item(119, 3).Value = "sjer5xc4"
$ws.Cells.Item(120, 3).Value = "sjer5xc5"
$ws.Cells.Item(121, 3).Value = "sjer5xc6"

# ---- "sjer20xc2" was entered ahead of the "sjer10x" block -----------------
$ws.Cells.Item(129, 3).Value = "sjer20xc2"

# ---- then the rest of the "sjer10x" block ----------------------------------
$ws.Cells.Item(123, 3).Value = "sjer10xc2"
$ws.Cells.Item(124, 3).Value = "sjer10xc3"
$ws.Cells.Item(125, 3).Value = "sjer10xc4"
$ws.Cells.Item(126, 3).Value = "sjer10xc5"
$ws.Cells.Item(127, 3).Value = "sjer10xc6"

# ---- then the rest of the "sjer20x" block ----------------------------------
$ws.Cells.Item(130, 3).Value = "sjer20xc3"
$ws.Cells.Item(131, 3).Value = "sjer20xc4"
$ws.Cells.Item(132, 3).Value = "sjer20xc5"
$ws.Cells.Item(133, 3).Value = "sjer20xc6"

# ---- Column D (Dilution Factor) --------------------------------------------
110..115 | ForEach-Object { $ws.Cells.Item($_, 4).Value = 1 }
116..121 | ForEach-Object { $ws.Cells.Item($_, 4).Value = 5 }
122..127 | ForEach-Object { $ws.Cells.Item($_, 4).Value = 10 }
128..133 | ForEach-Object { $ws.Cells.Item($_, 4).Value = 20 }

# ---- Restore the view: scrolled so row 99 is at the top, with the new
#      dilution-factor values (D110:D133) selected -------------------------
$ws.Range("A99").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 99
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D110:D133").Select() | Out-Null
